# Added code for New Word Mail Merge Template functionality
$wb = $excel.ActiveWorkbook

# --- New "Docs" worksheet, appended after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$docs = $wb.Worksheets.Add($null, $lastSheet)
$docs.Name = "Docs"

# Header row (yellow fill, matching the other sheets' header style)
$headers = @("title", "description", "version", "contact", "client", "prospect", "task", "case", "tags")
for ($col = 1; $col -le $headers.Length; $col++) {
    $docs.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$docs.Range($docs.Cells.Item(1, 1), $docs.Cells.Item(1, $headers.Length)).Interior.Color = 65535

# Data rows - version numbers are entered as text (leading apostrophe),
# matching how the other sheets store numeric-looking identifiers as text.
$row2 = @("Doc1", "Document 1", "'0.1", "abcd", "efgh", "ijkl", "mnop", "qrst", "uvwxyz")
for ($col = 1; $col -le $row2.Length; $col++) {
    $docs.Cells.Item(2, $col).Value = $row2[$col - 1]
}

$row3 = @("Doc2", "Document 2", "'1.9", "dcba", "hgfe", "lkji", "ponm", "tsrq", "zyxwvy")
for ($col = 1; $col -le $row3.Length; $col++) {
    $docs.Cells.Item(3, $col).Value = $row3[$col - 1]
}

# Column B (description) is wide enough that it needs to be auto-fit
$docs.Columns.Item(2).AutoFit()

# --- Previously active "Calls" sheet now selects the header row instead ---
$calls = $wb.Worksheets.Item("Calls")
$calls.Rows.Item(1).Select() | Out-Null

# Selection left at J3 on the new sheet, which becomes the active tab
$docs.Activate() | Out-Null
$docs.Cells.Item(3, 10).Select() | Out-Null
